$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 279, pushing existing rows 279:392 down to 280:393
$ws.Rows.Item(279).Insert()

# Populate the newly inserted row 279 with the new weekly record
$ws.Range("A279").Value = 3
$ws.Range("B279").Value = "Femacal de La Calera"
$ws.Range("C279").Value = "Coquimbo"
$ws.Range("D279").Value = 44924
$ws.Range("E279").Value = 5
$ws.Range("F279").Value = 100112039
$ws.Range("G279").Value = "Ciboulette"
$ws.Range("H279").Value = "Sin especificar"
$ws.Range("I279").Value = "Primera"
$ws.Range("J279").Value = 260
$ws.Range("K279").Value = 1500
$ws.Range("L279").Value = 1800
$ws.Range("M279").Value = 1673
$ws.Range("N279").Value = "$/docena de atados"
$ws.Range("O279").Value = "Provincia de Quillota"
$ws.Range("P279").Value = 558
$ws.Range("Q279").Value = 3
$ws.Range("R279").Value = "Hortaliza"
